# gsc-export update: append the next day's breadcrumb row to the Chart sheet.
#
# The GSC export tool appends one row per day. "2025-11-11" was the most
# recent date in the workbook; this adds "2025-11-12" with 0 invalid /
# 43 valid items.
#
# Note: writing a date-shaped literal straight into Range.Value lets Excel's
# normal typed-input behaviour kick in (the string gets parsed into a date
# serial + a date number format), but every other date in column A on this
# sheet is stored as plain text. So we pre-format the cell as Text ("@")
# before assigning the value to keep it a literal string, then ClearFormats
# so the cell's style matches the rest of the (unformatted, default-style)
# column instead of picking up the "@" text format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$row = 39

$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-11-12"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 43
